# Forced risk css class to lower case
# The "Risk" column (J) previously stored a numeric risk-level code
# (1/2/3). Replace those numbers with their lower-case textual labels.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = "medium risk"
$ws.Range("J3").Value = "low risk"
$ws.Range("J4").Value = "high risk"
$ws.Range("J5").Value = "high risk"
$ws.Range("J6").Value = "low risk"

# Match the author's final active-cell selection.
$ws.Range("J6").Select()
